# Fixed update to excel issue
# - Rename "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
# - Rename "Requested quantity" header on "Monthly Trend"  -> "Monthly_PO_Qty"
# - Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- Rename the B1 headers on the two existing sheets -----------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the formatting already present in the workbook (bold header style,
# date number-format on column A) instead of inventing new styles.
$wsMonthly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsMonthly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)  # xlPasteFormats

# --- Header row ---------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows (row 2 .. row 17) ----------------------------------------
$forecastData = @(
    @(44934.99999999999, 8, 4.918481137685623,  11.57186263128267),
    @(44941.99999999999, 8, 5.360041165015196,  10.91997128414836),
    @(44948.99999999999, 8, 4.633927869753209,  11.41755755339269),
    @(45025.99999999999, 8, 4.493767120341037,  10.72859397382515),
    @(45060.99999999999, 7, 4.098561918321864,  10.43629589201068),
    @(45067.99999999999, 7, 4.033186942149498,  10.3320153700349),
    @(45144.99999999999, 7, 3.674128886466983,  9.918718997156779),
    @(45151.99999999999, 7, 3.356949808334104,  9.615860772880975),
    @(45158.99999999999, 7, 3.306775158280113,  9.689551792101224),
    @(45165.99999999999, 7, 3.561931264353853,  9.803118911705949),
    @(45172.99999999999, 6, 3.278115073580005,  9.595959166032456),
    @(45179.99999999999, 6, 3.322773931063194,  9.588497216015824),
    @(45186.99999999999, 6, 3.061838374211193,  9.766940665463759),
    @(45193.99999999999, 6, 3.189974884151522,  9.604577532380707),
    @(45200.99999999999, 6, 2.973741688550744,  9.550023524437442),
    @(45207.99999999999, 6, 3.083662023533822,  9.105711609681569)
)

$row = 2
foreach ($record in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $record[0]
    $wsForecast.Cells.Item($row, 2).Value = $record[1]
    $wsForecast.Cells.Item($row, 3).Value = $record[2]
    $wsForecast.Cells.Item($row, 4).Value = $record[3]
    $row++
}

$wsForecast.Range("A1").Select() | Out-Null
